$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list with latest price/volume(1h) snapshot.
# D-column price cells are forced to Text format first since many values
# look numeric (e.g. "2.18", "0.998") and Excel would otherwise coerce them
# into real numbers, losing the original text formatting (e.g. trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '52.104.10'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.997.08'
$ws.Range("E3").Value = '  +2.64%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.46'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.55'
$ws.Range("E6").Value = '  -2.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.557'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("E9").Value = '  -3.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.07'
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("E12").Value = '  -4.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.00'
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.468.63'
$ws.Range("E14").Value = '  +2.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.60'
$ws.Range("E15").Value = '  -4.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.995.45'
$ws.Range("E16").Value = '  +2.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.01'
$ws.Range("E17").Value = '  +3.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '52.141.88'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.42'
$ws.Range("E19").Value = '  +5.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.47'
$ws.Range("E20").Value = '  -1.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.49'
$ws.Range("E21").Value = '  -4.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0968'
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.02'
$ws.Range("E23").Value = '  -2.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '263.81'
$ws.Range("E24").Value = '  -2.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.71'
$ws.Range("E25").Value = '  -2.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.176'
$ws.Range("E26").Value = '  -3.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.79'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.42'
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.38'
$ws.Range("E31").Value = '  +4.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.21'
$ws.Range("E32").Value = '  -3.50%  '
$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.18'
$ws.Range("E33").Value = '  +15.63%  '
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '35.77'
$ws.Range("E34").Value = '  -4.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '51.07'
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0438'
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.31'
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.83'
$ws.Range("E39").Value = '  +3.22%  '
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.58'
$ws.Range("E41").Value = '  -3.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.116'
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.89'
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '124.37'
$ws.Range("E44").Value = '  +7.81%  '
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.117.20'
$ws.Range("E46").Value = '  -0.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  -3.89%  '
$ws.Range("E48").Value = '  -6.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.295.57'
$ws.Range("E49").Value = '  +2.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.242'
$ws.Range("E50").Value = '  -2.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0330'
$ws.Range("E51").Value = '  -0.16%  '
